$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B17").Value = "1"
$ws.Range("B22").Value = "16"
$ws.Range("B23").Value = "16"
$ws.Range("B24").Value = "16"

$ws.Range("B25").Select()
